$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old used range (A1:C9) contents before rewriting as the new layout is smaller,
# keeping existing formatting (e.g. the bold/bordered header style) intact.
$ws.Range("A1:C9").ClearContents()

# Header row - column D is brand new, so copy the existing header formatting
# (bold font, border, centered alignment) from C1 onto it first.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("A1").Value = "geral_modalidade"
$ws.Range("B1").Value = "total"
$ws.Range("C1").Value = "total_sucesso"
$ws.Range("D1").Value = "total_falha"

# Data rows
$ws.Range("A2").Value = "aon"
$ws.Range("B2").Value = 1335
$ws.Range("C2").Value = 830
$ws.Range("D2").Value = 505

$ws.Range("A3").Value = "flex"
$ws.Range("B3").Value = 1468
$ws.Range("C3").Value = 1383
$ws.Range("D3").Value = 85

$ws.Range("A4").Value = "sub"
$ws.Range("B4").Value = 684
$ws.Range("C4").Value = 152
$ws.Range("D4").Value = 532
